$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update cost/unit-cost totals ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 1207.81123125
$schedule.Range("F2").Value = 19.97042379712302

# --- Sheet "Detailed": update price values / type labels for run 186 ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B8").Value = 56.98

$detailed.Range("B12").Value = 57.06

$detailed.Range("B13").Value = 65
$detailed.Range("C13").Value = "historical"

$detailed.Range("C14").Value = "historical"

$detailed.Range("B16").Value = 23.21022
$detailed.Range("B17").Value = 9.82175
$detailed.Range("B18").Value = 0.7
$detailed.Range("B19").Value = 22.07
$detailed.Range("B20").Value = 0.02916
$detailed.Range("B21").Value = 0.7
$detailed.Range("B22").Value = 4.00709

$detailed.Range("B25").Value = 35.26454

$detailed.Range("B33").Value = 4.22801
$detailed.Range("B34").Value = 6.31933
$detailed.Range("B35").Value = -5.03039
$detailed.Range("B36").Value = -3.04696
$detailed.Range("B37").Value = 6.00385
$detailed.Range("B38").Value = 44.57754
$detailed.Range("B39").Value = 56.98

$detailed.Range("B41").Value = 59.29218
$detailed.Range("B42").Value = 59.36731

$wb.Save()
